$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (price) contain numeric-looking text that must remain text,
# matching the original inlineStr cell type. Force Text format first so
# Excel does not auto-convert these into real numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.586.88"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.754.60"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "594.34"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "167.04"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("D7").Value = "3.752.07"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("D11").Value = "6.46"
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  -7.50%  "
$ws.Range("D14").Value = "36.10"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").Value = "4.386.78"
$ws.Range("D16").Value = "3.759.02"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "68.568.86"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "17.98"
$ws.Range("E18").Value = "  -4.92%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").Value = "  -3.16%  "
$ws.Range("D21").Value = "10.76"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").Value = "464.82"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").Value = "0.698"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").Value = "84.19"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").Value = "0.0000146"
$ws.Range("E25").Value = "  -3.51%  "
$ws.Range("E26").Value = "  -3.06%  "
$ws.Range("D27").Value = "11.95"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -4.75%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "3.902.94"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("E31").Value = "  -5.35%  "
$ws.Range("D32").Value = "7.31"
$ws.Range("E32").Value = "  -3.92%  "
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("E34").Value = "  -3.54%  "
$ws.Range("D35").Value = "9.19"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D37").Value = "3.708.29"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("E38").Value = "  -4.03%  "
$ws.Range("D39").Value = "3.40"
$ws.Range("E39").Value = "  -9.95%  "
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").Value = "5.79"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").Value = "44.05"
$ws.Range("E45").Value = "  +9.64%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.303"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "1.92"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "8.49"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D49").Value = "46.57"
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("D50").Value = "145.89"
$ws.Range("E50").Value = "  +3.52%  "
$ws.Range("D51").Value = "387.72"
$ws.Range("E51").Value = "  -4.29%  "

# Restore the default (unstyled) cell style on the price column now that
# the text values are locked in, so no stray formatting is left behind.
$ws.Range("D2:D51").Style = "Normal"
